$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = $origStyle
}

Set-TextValue 'D2' '69.223.11'
Set-TextValue 'E2' '  +1.69%  '
Set-TextValue 'D3' '3.386.96'
Set-TextValue 'E3' '  +1.36%  '
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.06%  '
Set-TextValue 'D5' '587.58'
Set-TextValue 'E5' '  +0.92%  '
Set-TextValue 'D6' '180.24'
Set-TextValue 'E6' '  +1.68%  '
Set-TextValue 'E7' '  -0.03%  '
Set-TextValue 'E8' '  +0.97%  '
Set-TextValue 'E9' '  +6.60%  '
Set-TextValue 'D10' '0.591'
Set-TextValue 'E10' '  +1.57%  '
Set-TextValue 'D11' '48.60'
Set-TextValue 'E11' '  +3.73%  '
Set-TextValue 'E12' '  +2.72%  '
Set-TextValue 'D13' '676.72'
Set-TextValue 'E13' '  -2.00%  '
Set-TextValue 'D14' '8.64'
Set-TextValue 'E14' '  +2.13%  '
Set-TextValue 'D15' '3.929.69'
Set-TextValue 'E15' '  +1.26%  '
Set-TextValue 'D16' '69.248.76'
Set-TextValue 'E16' '  +1.70%  '
Set-TextValue 'B17' 'TRON'
Set-TextValue 'C17' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D17' '0.120'
Set-TextValue 'E17' '  +1.66%  '
Set-TextValue 'B18' 'WrappedEther'
Set-TextValue 'C18' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D18' '3.382.56'
Set-TextValue 'E18' '  +1.30%  '
Set-TextValue 'D19' '17.68'
Set-TextValue 'E19' '  +1.39%  '
Set-TextValue 'D20' '11.33'
Set-TextValue 'E20' '  +2.26%  '
Set-TextValue 'D21' '0.903'
Set-TextValue 'E21' '  +0.44%  '
Set-TextValue 'E22' '  +0.99%  '
Set-TextValue 'D23' '17.10'
Set-TextValue 'E23' '  +0.46%  '
Set-TextValue 'D24' '103.36'
Set-TextValue 'E24' '  +4.03%  '
Set-TextValue 'E25' '  +0.69%  '
Set-TextValue 'E26' '  +0.76%  '
Set-TextValue 'D27' '9.60'
Set-TextValue 'E27' '  +0.44%  '
Set-TextValue 'D28' '34.16'
Set-TextValue 'E28' '  +3.29%  '
Set-TextValue 'D29' '8.72'
Set-TextValue 'E29' '  +1.55%  '
Set-TextValue 'D30' '7.00'
Set-TextValue 'E30' '  -1.30%  '
Set-TextValue 'E31' '  +1.60%  '
Set-TextValue 'E32' '  +10.60%  '
Set-TextValue 'D33' '554.81'
Set-TextValue 'E33' '  -3.00%  '
Set-TextValue 'E34' '  +0.54%  '
Set-TextValue 'D35' '57.99'
Set-TextValue 'E35' '  +1.06%  '
Set-TextValue 'E36' '  -0.05%  '
Set-TextValue 'D37' '3.687.36'
Set-TextValue 'E37' '  -0.97%  '
Set-TextValue 'E38' '  +6.59%  '
Set-TextValue 'D39' '35.02'
Set-TextValue 'E39' '  +1.02%  '
Set-TextValue 'E40' '  +1.25%  '
Set-TextValue 'B41' 'Fetch.AI'
Set-TextValue 'C41' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D41' '2.68'
Set-TextValue 'E41' '  +0.49%  '
Set-TextValue 'B42' 'PEPE'
Set-TextValue 'C42' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D42' '0.0₃0702'
Set-TextValue 'E42' '  +3.75%  '
Set-TextValue 'E43' '  +0.38%  '
Set-TextValue 'D44' '0.0423'
Set-TextValue 'E44' '  +3.72%  '
Set-TextValue 'D45' '3.28'
Set-TextValue 'E45' '  -1.48%  '
Set-TextValue 'E46' '  +0.08%  '
Set-TextValue 'E47' '  +0.70%  '
Set-TextValue 'E48' '  +5.45%  '
Set-TextValue 'E49' '  -0.07%  '
Set-TextValue 'D50' '131.82'
Set-TextValue 'E50' '  +1.50%  '
Set-TextValue 'D51' '2.58'
Set-TextValue 'E51' '  +1.30%  '
